$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - use Text number format to preserve exact string
# representation (avoids Excel auto-converting numeric-looking strings to numbers,
# which would drop trailing zeros / change formatting), then restore the original
# cell style so no formatting changes leak into the saved file.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D8",
    "D12",
    "D13",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D30",
    "D31",
    "D32",
    "D34",
    "D38",
    "D39",
    "D41",
    "D42",
    "D44",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
$priceValues = @(
    "63.670.90",
    "2.475.64",
    "577.74",
    "149.03",
    "0.542",
    "5.31",
    "27.27",
    "63.502.92",
    "2.482.84",
    "11.56",
    "7.36",
    "329.20",
    "1.97",
    "1.00",
    "67.50",
    "634.87",
    "8.94",
    "2.596.88",
    "8.42",
    "0.997",
    "0.147",
    "5.20",
    "5.55",
    "19.03",
    "146.59",
    "2.70",
    "150.81",
    "21.18",
    "0.612",
    "0.0240",
    "0.0929",
    "0.749"
)

for ($i = 0; $i -lt $priceCells.Length; $i++) {
    $rng = $ws.Range($priceCells[$i])
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $priceValues[$i]
    $rng.Style = $origStyle
}

# Volume(1h) percentage (column E) updates - plain text assignment is safe since
# these strings contain "%" and surrounding spaces, so Excel will never interpret
# them as numeric values.
$volCells = @(
    "E2",
    "E3",
    "E4",
    "E5",
    "E6",
    "E8",
    "E9",
    "E11",
    "E12",
    "E13",
    "E14",
    "E16",
    "E17",
    "E18",
    "E19",
    "E20",
    "E21",
    "E22",
    "E23",
    "E24",
    "E25",
    "E26",
    "E27",
    "E28",
    "E29",
    "E30",
    "E31",
    "E32",
    "E33",
    "E34",
    "E35",
    "E36",
    "E37",
    "E38",
    "E39",
    "E40",
    "E41",
    "E42",
    "E43",
    "E44",
    "E45",
    "E46",
    "E47",
    "E48",
    "E49",
    "E50",
    "E51"
)
$volValues = @(
    "  +2.66%  ",
    "  +2.20%  ",
    "  +0.06%  ",
    "  +2.64%  ",
    "  +4.07%  ",
    "  +2.02%  ",
    "  +4.44%  ",
    "  +3.97%  ",
    "  +2.21%  ",
    "  +4.24%  ",
    "  +6.24%  ",
    "  +2.50%  ",
    "  +2.38%  ",
    "  +2.15%  ",
    "  +7.72%  ",
    "  +2.77%  ",
    "  +1.67%  ",
    "  +13.17%  ",
    "  +0.03%  ",
    "  +0.63%  ",
    "  +14.10%  ",
    "  +1.94%  ",
    "  +14.19%  ",
    "  +2.10%  ",
    "  +9.81%  ",
    "  +2.57%  ",
    "  -0.11%  ",
    "  -0.52%  ",
    "  +3.67%  ",
    "  +9.52%  ",
    "  +3.37%  ",
    "  -0.19%  ",
    "  +2.25%  ",
    "  +1.77%  ",
    "  +2.00%  ",
    "  +2.47%  ",
    "  -3.98%  ",
    "  +20.92%  ",
    "  +0.56%  ",
    "  +2.27%  ",
    "  +4.09%  ",
    "  +4.51%  ",
    "  +6.90%  ",
    "  +2.52%  ",
    "  +5.36%  ",
    "  +0.95%  ",
    "  +4.53%  "
)

for ($i = 0; $i -lt $volCells.Length; $i++) {
    $ws.Range($volCells[$i]).Value = $volValues[$i]
}

Write-Host "Applied $($priceCells.Length) price updates and $($volCells.Length) volume updates."
